# dodelani vizualizace zobrazeni logu + oprava logiky zapisovani
#
# Fills in the task_settings sheet with 3 task-log rows (A:G) and gives the
# log column (G) a monospace "console" look (Consolas, orange-ish text,
# vertically centered). Also keeps the time-formatted helper cell (E1) as
# plain text, fixes the sheet's page setup and leaves the G1 cell selected
# like the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task_settings")

# Helper: write a value that LOOKS like a plain number ("292", "50", ...)
# while keeping it stored as text, exactly like the source log lines. We
# briefly force a Text number format so Excel does not coerce the literal
# into a numeric cell, then reset the style back to Normal so no stray
# per-cell format lingers on these otherwise unstyled cells.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 1 - TRIMAZKON_task_3
$ws.Cells.Item(1,1).Value = 'TRIMAZKON_task_3'
$ws.Cells.Item(1,2).Value = 'C:/Users/jakub.hlavacek.local/Desktop/JHV/test_images/Keyence/_503_Witte/datumovka/A/Height_test/'
Set-TextValue $ws.Cells.Item(1,3) '292'
Set-TextValue $ws.Cells.Item(1,4) '50'
$ws.Cells.Item(1,5).Value = '12:00'
$ws.Cells.Item(1,6).Value = '17.12.2024 16:11:56'
$ws.Cells.Item(1,7).Value = '|||Datum provedení: 17.12.2024 16:12:19||Zkontrolováno: 161 souborů||Starších:      153 souborů||Smazáno:       103 souborů'

# Row 2 - TRIMAZKON_task_2
$ws.Cells.Item(2,1).Value = 'TRIMAZKON_task_2'
$ws.Cells.Item(2,2).Value = 'C:/Users/jakub.hlavacek.local/Desktop/JHV/test_images/Keyence/_503_Witte/datumovka/A/Height_test/'
Set-TextValue $ws.Cells.Item(2,3) '292'
Set-TextValue $ws.Cells.Item(2,4) '998'
$ws.Cells.Item(2,5).Value = '12:00'
$ws.Cells.Item(2,6).Value = '17.12.2024 10:12:10'
$ws.Cells.Item(2,7).Value = '|||Datum: 17.122222210:12:26||Zkontrolováno: 161 souborů||Starších: 153 souborů||Smazáno: 0 souborů|||Datum: 17.12.2024 10:12:26||Zkontrolováno: 161 souborů||Starších: 153 souborů||Smazáno: 0 souborů'

# Row 3 - TRIMAZKON_task_1
$ws.Cells.Item(3,1).Value = 'TRIMAZKON_task_1'
$ws.Cells.Item(3,2).Value = 'C:/Users/jakub.hlavacek.local/Desktop/JHV/test_images/Keyence/_503_Witte/datumovka/A/Height_test/'
Set-TextValue $ws.Cells.Item(3,3) '100'
Set-TextValue $ws.Cells.Item(3,4) '200'
$ws.Cells.Item(3,5).Value = '6:00'
$ws.Cells.Item(3,6).Value = '17.12.2024 16:03:53'
$ws.Cells.Item(3,7).Value = '|||Datum provedení: 17.12.2024 16:04:58||Zkontrolováno: 161 souborů||Starších:      153 souborů||Smazáno:       0 souborů|||Datum provedení: 17.12.2024 16:05:16||Zkontrolováno: 161 souborů||Starších:      153 souborů||Smazáno:       0 souborů|||Datum provedení: 17.12.2024 16:08:31||Zkontrolováno: 161 souborů||Starších:      153 souborů||Smazáno:       0 souborů|||Datum provedení: 17.12.2024 16:17:29||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů'

# Log column look & feel: Consolas, orange-ish (CE9178), vertically centered.
$logRange = $ws.Range("G1:G2")
$logRange.Font.Name = "Consolas"
$logRange.Font.Color = 7901646
$logRange.VerticalAlignment = -4108

# Page setup (portrait / A4) for printing this sheet.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# Leave the selection on G1, matching the authored state.
$ws.Activate() | Out-Null
$ws.Range("G1").Select() | Out-Null
